# Apply "Research info" column (O) content for the Details sheet table (Table2).
# Values are entered in the same order the original author typed them so that the
# workbook's shared-string table gets populated in a matching sequence (dedup on repeats).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Details")

$ws.Range("O3").Value = "User is travelling out of country, personal choice, doesn,t want misuse. Can be anything. It's a user choice. That is the use case. 
Check if OTP authentication is allowed or not- (if a Mobile Numberis associated with multiple UIN/VID then the country can term it as Weak based on the number)"
$ws.Range("O5").Value = "Check if OTP authentication is allowed or not- (if a Mobile Numberis associated with multiple UIN/VID then the country can term it as Weak based on the number)"
$ws.Range("O4").Value = "User likes multiple copies, use lost old copy. Can be anything. Provision to keep track on no of reprint required, country can reject free request after a X upper limit, X can be zero, country to allow payment based reprint, thus integration point with payment gateway will be required for SI customisation etc. 
validation and the Interface for Payment gateway will be part of Resident portal(business Logic)- by SI
The service should cater to any service received for Re-printing
Check if OTP authentication is allowed or not- (if a Mobile Numberis associated with multiple UIN/VID then the country can term it as Weak based on the number)"
$ws.Range("O7").Value = "Can update any demo graphic information, validations should be done by Portal.
Should cater to updation of demographic details as done by Reg. client
Check if OTP authentication is allowed or not- (if a Mobile Numberis associated with multiple UIN/VID then the country can term it as Weak based on the number)"
$ws.Range("O8").Value = "think and come back on what all status can the request have and we can rationalise. 
Check if OTP authentication is allowed or not- (if a Mobile Numberis associated with multiple UIN/VID then the country can term it as Weak based on the number)"
$ws.Range("O10").Value = "YES. Correct understanding. 
Check if OTP authentication is allowed or not- (if a Mobile Numberis associated with multiple UIN/VID then the country can term it as Weak based on the number)"
$ws.Range("O6").Value = "Check if OTP authentication is allowed or not- (if a Mobile Numberis associated with multiple UIN/VID then the country can term it as Weak based on the number)"
$ws.Range("O9").Value = "think and come back on what all status can the request have and we can rationalise. 
Check if OTP authentication is allowed or not- (if a Mobile Numberis associated with multiple UIN/VID then the country can term it as Weak based on the number)"
$ws.Range("O11").Value = "Check if OTP authentication is allowed or not- (if a Mobile Numberis associated with multiple UIN/VID then the country can term it as Weak based on the number)"

# The base "Check if OTP..." note (first entered in O5) is left-aligned by default;
# switch it to the default/general horizontal alignment to match the final formatting.
$ws.Range("O5").HorizontalAlignment = 1

# Restore view state: the author ended up zoomed in on column O around row 11.
$ws.Activate()
$excel.ActiveWindow.Zoom = 110
$ws.Range("O11").Select()

